$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value2 = "Datos actualizados a 27 de Mayo de 2020 a las 17:35"

# Update province/city names and case statistics (B=Casos totales, C=Casos activos, D=Recuperados, E=Muertes)
# Row 10: 'Bizkaia/Vizcaya' -> 'Bizkaia/Vizcaya'
$ws.Range("C10").Value2 = 0
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Value2 = 10332

# Row 12: 'Ciudad Real' -> 'Ciudad Real'
$ws.Range("C12").Value2 = 0
$ws.Range("D12").Value2 = 0
$ws.Range("E12").Value2 = 6464

# Row 13: 'Valencia/Valencia' -> 'Valencia/Valencia'
$ws.Range("C13").Value2 = 0
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 5609

# Row 15: 'Zaragoza' -> 'Zaragoza'
$ws.Range("C15").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 5287

# Row 17: 'Araba/Alava' -> 'Araba/Alava'
$ws.Range("C17").Value2 = 0
$ws.Range("D17").Value2 = 0
$ws.Range("E17").Value2 = 4868

# Row 18: 'Valladolid' -> 'Valladolid'
$ws.Range("C18").Value2 = 0
$ws.Range("D18").Value2 = 0
$ws.Range("E18").Value2 = 4393

# Row 19: 'Salamanca' -> 'Salamanca'
$ws.Range("C19").Value2 = 0
$ws.Range("D19").Value2 = 0
$ws.Range("E19").Value2 = 4152

# Row 21: 'Toledo' -> 'Toledo'
$ws.Range("C21").Value2 = 0
$ws.Range("D21").Value2 = 0
$ws.Range("E21").Value2 = 3872

# Row 22: 'Alacant/Alicante' -> 'Alacant/Alicante'
$ws.Range("C22").Value2 = 0
$ws.Range("D22").Value2 = 0
$ws.Range("E22").Value2 = 3794

# Row 23: 'Albacete' -> 'Albacete'
$ws.Range("C23").Value2 = 0
$ws.Range("D23").Value2 = 0
$ws.Range("E23").Value2 = 3775

# Row 24: 'Leon' -> 'Leon'
$ws.Range("C24").Value2 = 0
$ws.Range("D24").Value2 = 0
$ws.Range("E24").Value2 = 3569

# Row 25: 'Segovia' -> 'Segovia'
$ws.Range("C25").Value2 = 0
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 3413

# Row 26: 'Gipuzkoa/Guipuzcoa' -> 'Gipuzkoa/Guipuzcoa'
$ws.Range("C26").Value2 = 0
$ws.Range("D26").Value2 = 0
$ws.Range("E26").Value2 = 3116

# Row 28: 'Malaga' -> 'Malaga'
$ws.Range("C28").Value2 = 0
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 2758

# Row 29: 'Burgos' -> 'Burgos'
$ws.Range("C29").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 2746

# Row 30: 'Sevilla' -> 'Sevilla'
$ws.Range("C30").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 2423

# Row 31: 'Granada' -> 'Granada'
$ws.Range("C31").Value2 = 0
$ws.Range("E31").Value2 = 2413

# Row 33: 'Gran Canaria' -> 'Soria'
$ws.Range("A33").Value2 = "Soria"
$ws.Range("B33").Value2 = 2290
$ws.Range("C33").Value2 = 0
$ws.Range("D33").Value2 = 0
$ws.Range("E33").Value2 = 2290

# Row 34: 'Soria' -> 'Tenerife'
$ws.Range("A34").Value2 = "Tenerife"
$ws.Range("B34").Value2 = 2280
$ws.Range("C34").Value2 = 0
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 2280

# Row 35: 'Tenerife' -> 'Cantabria'
$ws.Range("A35").Value2 = "Cantabria"
$ws.Range("B35").Value2 = 2246
$ws.Range("C35").Value2 = 1981
$ws.Range("D35").Value2 = 62
$ws.Range("E35").Value2 = 203

# Row 36: 'Cantabria' -> 'Caceres'
$ws.Range("A36").Value2 = "Caceres"
$ws.Range("B36").Value2 = 1973
$ws.Range("C36").Value2 = 0
$ws.Range("D36").Value2 = 0
$ws.Range("E36").Value2 = 1973

# Row 37: 'Caceres' -> 'A Coruña'
$ws.Range("A37").Value2 = "A Coruña"
$ws.Range("B37").Value2 = 1969
$ws.Range("C37").Value2 = 333
$ws.Range("D37").Value2 = 1788
$ws.Range("E37").Value2 = 67

# Row 38: 'A Coruña' -> 'Avila'
$ws.Range("A38").Value2 = "Avila"
$ws.Range("B38").Value2 = 1935
$ws.Range("C38").Value2 = 0
$ws.Range("D38").Value2 = 0
$ws.Range("E38").Value2 = 1935

# Row 39: 'Avila' -> 'Murcia'
$ws.Range("A39").Value2 = "Murcia"
$ws.Range("B39").Value2 = 1587
$ws.Range("C39").Value2 = 2180
$ws.Range("D39").Value2 = 0
$ws.Range("E39").Value2 = 148

# Row 40: 'Murcia' -> 'Pontevedra'
$ws.Range("A40").Value2 = "Pontevedra"
$ws.Range("B40").Value2 = 1536
$ws.Range("C40").Value2 = 333
$ws.Range("D40").Value2 = 1411
$ws.Range("E40").Value2 = 30

# Row 41: 'Pontevedra' -> 'Castello/Castellon'
$ws.Range("A41").Value2 = "Castello/Castellon"
$ws.Range("B41").Value2 = 1486
$ws.Range("C41").Value2 = 0
$ws.Range("D41").Value2 = 0
$ws.Range("E41").Value2 = 1486

# Row 42: 'Castello/Castellon' -> 'Jaen'
$ws.Range("A42").Value2 = "Jaen"
$ws.Range("B42").Value2 = 1387
$ws.Range("C42").Value2 = 0
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 1387

# Row 43: 'Jaen' -> 'Cordoba'
$ws.Range("A43").Value2 = "Cordoba"
$ws.Range("B43").Value2 = 1331
$ws.Range("C43").Value2 = 0
$ws.Range("D43").Value2 = 0
$ws.Range("E43").Value2 = 1331

# Row 44: 'Cordoba' -> 'Guadalajara'
$ws.Range("A44").Value2 = "Guadalajara"
$ws.Range("B44").Value2 = 1266
$ws.Range("C44").Value2 = 0
$ws.Range("E44").Value2 = 1266

# Row 45: 'Guadalajara' -> 'Cuenca'
$ws.Range("A45").Value2 = "Cuenca"
$ws.Range("B45").Value2 = 1241
$ws.Range("C45").Value2 = 0
$ws.Range("D45").Value2 = 0
$ws.Range("E45").Value2 = 1241

# Row 46: 'Cuenca' -> 'Cadiz'
$ws.Range("A46").Value2 = "Cadiz"
$ws.Range("B46").Value2 = 1240
$ws.Range("C46").Value2 = 0
$ws.Range("D46").Value2 = 0
$ws.Range("E46").Value2 = 1240

# Row 47: 'Cadiz' -> 'Palencia'
$ws.Range("A47").Value2 = "Palencia"
$ws.Range("B47").Value2 = 1205
$ws.Range("C47").Value2 = 0
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 1205

# Row 48: 'Palencia' -> 'Huesca'
$ws.Range("A48").Value2 = "Huesca"
$ws.Range("B48").Value2 = 1115
$ws.Range("C48").Value2 = 0
$ws.Range("D48").Value2 = 0
$ws.Range("E48").Value2 = 1115

# Row 49: 'Huesca' -> 'Zamora'
$ws.Range("A49").Value2 = "Zamora"
$ws.Range("B49").Value2 = 993
$ws.Range("C49").Value2 = 0
$ws.Range("D49").Value2 = 0
$ws.Range("E49").Value2 = 993

# Row 50: 'Zamora' -> 'Badajoz'
$ws.Range("A50").Value2 = "Badajoz"
$ws.Range("B50").Value2 = 962
$ws.Range("C50").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 962

# Row 51: 'Badajoz' -> 'Ourense'
$ws.Range("A51").Value2 = "Ourense"
$ws.Range("B51").Value2 = 751
$ws.Range("C51").Value2 = 333
$ws.Range("D51").Value2 = 660
$ws.Range("E51").Value2 = 22

# Row 52: 'Ourense' -> 'Teruel'
$ws.Range("A52").Value2 = "Teruel"
$ws.Range("B52").Value2 = 664
$ws.Range("C52").Value2 = 0
$ws.Range("D52").Value2 = 0
$ws.Range("E52").Value2 = 664

# Row 53: 'Teruel' -> 'Lugo'
$ws.Range("A53").Value2 = "Lugo"
$ws.Range("B53").Value2 = 586
$ws.Range("C53").Value2 = 333
$ws.Range("D53").Value2 = 520
$ws.Range("E53").Value2 = 11

# Row 54: 'Lugo' -> 'Gran Canaria'
$ws.Range("A54").Value2 = "Gran Canaria"
$ws.Range("B54").Value2 = 563
$ws.Range("C54").Value2 = 0
$ws.Range("D54").Value2 = 0
$ws.Range("E54").Value2 = 563

# Row 55: 'Almeria' -> 'Almeria'
$ws.Range("C55").Value2 = 0
$ws.Range("E55").Value2 = 498

# Row 56: 'Huelva' -> 'Huelva'
$ws.Range("C56").Value2 = 0
$ws.Range("E56").Value2 = 400

# Row 60: 'La Palma' -> 'La Palma'
$ws.Range("C60").Value2 = 0
$ws.Range("D60").Value2 = 0
$ws.Range("E60").Value2 = 95

# Row 61: 'Lanzarote' -> 'Lanzarote'
$ws.Range("C61").Value2 = 0
$ws.Range("D61").Value2 = 0
$ws.Range("E61").Value2 = 84

# Row 63: 'Fuerteventura' -> 'Fuerteventura'
$ws.Range("C63").Value2 = 0
$ws.Range("E63").Value2 = 23

# Row 66: 'La Gomera' -> 'La Gomera'
$ws.Range("C66").Value2 = 0
$ws.Range("E66").Value2 = 8

# Row 68: 'El Hierro' -> 'El Hierro'
$ws.Range("C68").Value2 = 0
$ws.Range("E68").Value2 = 3

